$d = $word.ActiveDocument

# Insert two new paragraphs after the "2022年6月3日星期五" paragraph
# (paragraph 5), before the "中雨，今天是农历五月初五，中国传统端午节。"
# paragraph, so they inherit the eastAsia-hinted paragraph formatting.
$anchor = $d.Paragraphs.Item(5)
$anchor.Range.InsertParagraphAfter()

$p1 = $d.Paragraphs.Item(6)
$p1.Range.Text = "中雨，今天是农历五月初五，中国传统端午节，这一天我们要吃粽子，赛龙舟。"
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(7)
$p2.Range.Text = "2022年6月7日星期二"

# Update the text of the original last-dated paragraph.
$d.Content.Find.Execute("中雨，今天是农历五月初五，中国传统端午节。", $true, $false, $false, $false, $false,
                         $true, 1, $false, "晴，今天是高考第一天，上午考语文，下午考数学。", 2)
